$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.011658668518066
$ws.Range("B1").Value = 2.130697011947632
$ws.Range("C1").Value = 5.698097229003906
$ws.Range("D1").Value = 0.9191138744354248
$ws.Range("E1").Value = 1.000214576721191
